$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "308.07"
$ws.Cells.Item(2, 4).NumberFormat = "General"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "0.84%"
$ws.Cells.Item(2, 5).NumberFormat = "General"
$ws.Cells.Item(2, 5).Style = "Normal"

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "38.86"
$ws.Cells.Item(3, 4).NumberFormat = "General"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = "8.71%"
$ws.Cells.Item(3, 5).NumberFormat = "General"
$ws.Cells.Item(3, 5).Style = "Normal"

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "5.102"
$ws.Cells.Item(4, 4).NumberFormat = "General"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = "1.06%"
$ws.Cells.Item(4, 5).NumberFormat = "General"
$ws.Cells.Item(4, 5).Style = "Normal"

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "0.08132"
$ws.Cells.Item(5, 4).NumberFormat = "General"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = "1.33%"
$ws.Cells.Item(5, 5).NumberFormat = "General"
$ws.Cells.Item(5, 5).Style = "Normal"

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "1.974"
$ws.Cells.Item(6, 4).NumberFormat = "General"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = "4.75%"
$ws.Cells.Item(6, 5).NumberFormat = "General"
$ws.Cells.Item(6, 5).Style = "Normal"

$ws.Cells.Item(7, 2).Value = "KuCoinToken"
$ws.Cells.Item(7, 3).Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "7.930"
$ws.Cells.Item(7, 4).NumberFormat = "General"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = "1.94%"
$ws.Cells.Item(7, 5).NumberFormat = "General"
$ws.Cells.Item(7, 5).Style = "Normal"

$ws.Cells.Item(8, 2).Value = "MXToken"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.9284"
$ws.Cells.Item(8, 4).NumberFormat = "General"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = "0.75%"
$ws.Cells.Item(8, 5).NumberFormat = "General"
$ws.Cells.Item(8, 5).Style = "Normal"

$ws.Cells.Item(9, 2).Value = "LiechtensteinCryptoassetsExchange"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.1444"
$ws.Cells.Item(9, 4).NumberFormat = "General"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = "13.12%"
$ws.Cells.Item(9, 5).NumberFormat = "General"
$ws.Cells.Item(9, 5).Style = "Normal"

$ws.Cells.Item(10, 2).Value = "WazirX"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.1957"
$ws.Cells.Item(10, 4).NumberFormat = "General"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = "2.52%"
$ws.Cells.Item(10, 5).NumberFormat = "General"
$ws.Cells.Item(10, 5).Style = "Normal"

$ws.Cells.Item(11, 2).Value = "MandalaExchangeToken"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.09096"
$ws.Cells.Item(11, 4).NumberFormat = "General"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = "0.27%"
$ws.Cells.Item(11, 5).NumberFormat = "General"
$ws.Cells.Item(11, 5).Style = "Normal"

$ws.Cells.Item(12, 2).Value = "BitrueCoin"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.03515"
$ws.Cells.Item(12, 4).NumberFormat = "General"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = "2.06%"
$ws.Cells.Item(12, 5).NumberFormat = "General"
$ws.Cells.Item(12, 5).Style = "Normal"

$ws.Cells.Item(13, 2).Value = "BitMartToken"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.09834"
$ws.Cells.Item(13, 4).NumberFormat = "General"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = "-0.19%"
$ws.Cells.Item(13, 5).NumberFormat = "General"
$ws.Cells.Item(13, 5).Style = "Normal"

$ws.Cells.Item(14, 2).Value = "BitForexToken"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.001402"
$ws.Cells.Item(14, 4).NumberFormat = "General"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = "-0.52%"
$ws.Cells.Item(14, 5).NumberFormat = "General"
$ws.Cells.Item(14, 5).Style = "Normal"

$ws.Cells.Item(15, 2).Value = "TigerCash"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.005796"
$ws.Cells.Item(15, 4).NumberFormat = "General"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = "-6.68%"
$ws.Cells.Item(15, 5).NumberFormat = "General"
$ws.Cells.Item(15, 5).Style = "Normal"

$ws.Cells.Item(16, 2).Value = "LEO"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "3.598"
$ws.Cells.Item(16, 4).NumberFormat = "General"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = "-4.19%"
$ws.Cells.Item(16, 5).NumberFormat = "General"
$ws.Cells.Item(16, 5).Style = "Normal"

$ws.Cells.Item(17, 2).Value = "GateToken"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "4.204"
$ws.Cells.Item(17, 4).NumberFormat = "General"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = "1.07%"
$ws.Cells.Item(17, 5).NumberFormat = "General"
$ws.Cells.Item(17, 5).Style = "Normal"

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "3.410"
$ws.Cells.Item(18, 4).NumberFormat = "General"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = "0.92%"
$ws.Cells.Item(18, 5).NumberFormat = "General"
$ws.Cells.Item(18, 5).Style = "Normal"

$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = "-0.02%"
$ws.Cells.Item(19, 5).NumberFormat = "General"
$ws.Cells.Item(19, 5).Style = "Normal"

$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = "-1.08%"
$ws.Cells.Item(20, 5).NumberFormat = "General"
$ws.Cells.Item(20, 5).Style = "Normal"

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "4.814"
$ws.Cells.Item(21, 4).NumberFormat = "General"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = "-6.87%"
$ws.Cells.Item(21, 5).NumberFormat = "General"
$ws.Cells.Item(21, 5).Style = "Normal"

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.2455"
$ws.Cells.Item(22, 4).NumberFormat = "General"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = "-5.76%"
$ws.Cells.Item(22, 5).NumberFormat = "General"
$ws.Cells.Item(22, 5).Style = "Normal"

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.04438"
$ws.Cells.Item(23, 4).NumberFormat = "General"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = "0.37%"
$ws.Cells.Item(23, 5).NumberFormat = "General"
$ws.Cells.Item(23, 5).Style = "Normal"

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "0.001228"
$ws.Cells.Item(24, 4).NumberFormat = "General"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = "-0.51%"
$ws.Cells.Item(24, 5).NumberFormat = "General"
$ws.Cells.Item(24, 5).Style = "Normal"

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.004847"
$ws.Cells.Item(25, 4).NumberFormat = "General"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = "4.96%"
$ws.Cells.Item(25, 5).NumberFormat = "General"
$ws.Cells.Item(25, 5).Style = "Normal"

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.0001303"
$ws.Cells.Item(27, 4).NumberFormat = "General"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = "4.10%"
$ws.Cells.Item(27, 5).NumberFormat = "General"
$ws.Cells.Item(27, 5).Style = "Normal"

$ws.Cells.Item(39, 5).NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = "8.49%"
$ws.Cells.Item(39, 5).NumberFormat = "General"
$ws.Cells.Item(39, 5).Style = "Normal"

$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = "-8.65%"
$ws.Cells.Item(40, 5).NumberFormat = "General"
$ws.Cells.Item(40, 5).Style = "Normal"

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.007465"
$ws.Cells.Item(41, 4).NumberFormat = "General"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = "-1.94%"
$ws.Cells.Item(41, 5).NumberFormat = "General"
$ws.Cells.Item(41, 5).Style = "Normal"

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.01012"
$ws.Cells.Item(42, 4).NumberFormat = "General"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = "-0.47%"
$ws.Cells.Item(42, 5).NumberFormat = "General"
$ws.Cells.Item(42, 5).Style = "Normal"

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.1362"
$ws.Cells.Item(43, 4).NumberFormat = "General"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = "0.63%"
$ws.Cells.Item(43, 5).NumberFormat = "General"
$ws.Cells.Item(43, 5).Style = "Normal"

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.002145"
$ws.Cells.Item(44, 4).NumberFormat = "General"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = "-0.36%"
$ws.Cells.Item(44, 5).NumberFormat = "General"
$ws.Cells.Item(44, 5).Style = "Normal"

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.01049"
$ws.Cells.Item(45, 4).NumberFormat = "General"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = "6.58%"
$ws.Cells.Item(45, 5).NumberFormat = "General"
$ws.Cells.Item(45, 5).Style = "Normal"

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.00006244"
$ws.Cells.Item(46, 4).NumberFormat = "General"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = "1.59%"
$ws.Cells.Item(46, 5).NumberFormat = "General"
$ws.Cells.Item(46, 5).Style = "Normal"

$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = "-0.01%"
$ws.Cells.Item(47, 5).NumberFormat = "General"
$ws.Cells.Item(47, 5).Style = "Normal"

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.003070"
$ws.Cells.Item(48, 4).NumberFormat = "General"
$ws.Cells.Item(48, 4).Style = "Normal"

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.001602"
$ws.Cells.Item(49, 4).NumberFormat = "General"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).NumberFormat = "@"
$ws.Cells.Item(49, 5).Value = "-3.51%"
$ws.Cells.Item(49, 5).NumberFormat = "General"
$ws.Cells.Item(49, 5).Style = "Normal"

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.00002103"
$ws.Cells.Item(50, 4).NumberFormat = "General"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).NumberFormat = "@"
$ws.Cells.Item(50, 5).Value = "-0.01%"
$ws.Cells.Item(50, 5).NumberFormat = "General"
$ws.Cells.Item(50, 5).Style = "Normal"

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.0002003"
$ws.Cells.Item(51, 4).NumberFormat = "General"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).NumberFormat = "@"
$ws.Cells.Item(51, 5).Value = "-0.01%"
$ws.Cells.Item(51, 5).NumberFormat = "General"
$ws.Cells.Item(51, 5).Style = "Normal"
